# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Primera / Segunda) for Coliflor at
# "Terminal La Palmera de La Serena", pushing the existing data down by
# two rows (Excel's standard "insert rows" shift-down behaviour).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 1195 (shifts 1195:1285 -> 1197:1287)
$ws.Range("A1195:A1196").EntireRow.Insert()

# New row 1195 - Calidad "Primera"
$ws.Cells.Item(1195, 1).Value = 8
$ws.Cells.Item(1195, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1195, 3).Value = "Coquimbo"
$ws.Cells.Item(1195, 4).Value = 45265
$ws.Cells.Item(1195, 5).Value = 4
$ws.Cells.Item(1195, 6).Value = 100112008
$ws.Cells.Item(1195, 7).Value = "Coliflor"
$ws.Cells.Item(1195, 8).Value = "Sin especificar"
$ws.Cells.Item(1195, 9).Value = "Primera"
$ws.Cells.Item(1195, 10).Value = 2000
$ws.Cells.Item(1195, 11).Value = 700
$ws.Cells.Item(1195, 12).Value = 800
$ws.Cells.Item(1195, 13).Value = 750
$ws.Cells.Item(1195, 14).Value = "$/unidad"
$ws.Cells.Item(1195, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(1195, 16).Value = 750
$ws.Cells.Item(1195, 17).Value = 1
$ws.Cells.Item(1195, 18).Value = "Hortaliza"

# New row 1196 - Calidad "Segunda"
$ws.Cells.Item(1196, 1).Value = 8
$ws.Cells.Item(1196, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1196, 3).Value = "Coquimbo"
$ws.Cells.Item(1196, 4).Value = 45265
$ws.Cells.Item(1196, 5).Value = 4
$ws.Cells.Item(1196, 6).Value = 100112008
$ws.Cells.Item(1196, 7).Value = "Coliflor"
$ws.Cells.Item(1196, 8).Value = "Sin especificar"
$ws.Cells.Item(1196, 9).Value = "Segunda"
$ws.Cells.Item(1196, 10).Value = 1100
$ws.Cells.Item(1196, 11).Value = 500
$ws.Cells.Item(1196, 12).Value = 600
$ws.Cells.Item(1196, 13).Value = 550
$ws.Cells.Item(1196, 14).Value = "$/unidad"
$ws.Cells.Item(1196, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(1196, 16).Value = 550
$ws.Cells.Item(1196, 17).Value = 1
$ws.Cells.Item(1196, 18).Value = "Hortaliza"
